$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so values like "1.018" are not
# auto-converted into numbers by Excel (locale decimal separator is ".").
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.828.36'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").Value = '1.855.56'
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("D4").Value = '1.018'
$ws.Range("E4").Value = '  -1.67%  '

$ws.Range("D5").Value = '320.33'
$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("D6").Value = '1.016'
$ws.Range("E6").Value = '  -1.75%  '

$ws.Range("D7").Value = '0.4358'
$ws.Range("E7").Value = '  -1.37%  '

$ws.Range("D8").Value = '0.3775'
$ws.Range("E8").Value = '  -0.83%  '

$ws.Range("D9").Value = '0.07413'
$ws.Range("E9").Value = '  -0.68%  '

$ws.Range("D10").Value = '0.8822'
$ws.Range("E10").Value = '  -0.17%  '

$ws.Range("D11").Value = '21.61'
$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").Value = '1.853.73'
$ws.Range("E12").Value = '  -0.97%  '

$ws.Range("D13").Value = '6.774'
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").Value = '5.484'
$ws.Range("E14").Value = '  -1.46%  '

$ws.Range("D15").Value = '0.07112'
$ws.Range("E15").Value = '  -1.59%  '

$ws.Range("D16").Value = '88.41'
$ws.Range("E16").Value = '  +5.52%  '

$ws.Range("D17").Value = '1.022'
$ws.Range("E17").Value = '  -1.68%  '

$ws.Range("D18").Value = '0.000009020'
$ws.Range("E18").Value = '  -0.92%  '

$ws.Range("D19").Value = '1.016'
$ws.Range("E19").Value = '  -1.72%  '

$ws.Range("D20").Value = '15.47'
$ws.Range("E20").Value = '  -0.42%  '

$ws.Range("D21").Value = '27.795.63'
$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D22").Value = '5.257'
$ws.Range("E22").Value = '  -1.07%  '

$ws.Range("D23").Value = '11.16'
$ws.Range("E23").Value = '  -2.56%  '

$ws.Range("D24").Value = '2.095.26'
$ws.Range("E24").Value = '  +0.37%  '

$ws.Range("D25").Value = '2.029'
$ws.Range("E25").Value = '  +4.80%  '

$ws.Range("D26").Value = '156.73'
$ws.Range("E26").Value = '  -1.27%  '

$ws.Range("D27").Value = '18.67'
$ws.Range("E27").Value = '  -1.09%  '

$ws.Range("D28").Value = '5.422'
$ws.Range("E28").Value = '  +1.54%  '

$ws.Range("D29").Value = '1.993'
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("D30").Value = '120.54'
$ws.Range("E30").Value = '  +2.35%  '

$ws.Range("D31").Value = '0.09031'
$ws.Range("E31").Value = '  -0.86%  '

$ws.Range("D32").Value = '1.228'
$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("D33").Value = '0.7685'
$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '4.557'
$ws.Range("E34").Value = '  -0.48%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.982'
$ws.Range("E35").Value = '  +2.60%  '

$ws.Range("D36").Value = '1.018'
$ws.Range("E36").Value = '  -1.73%  '

$ws.Range("D37").Value = '1.139'
$ws.Range("E37").Value = '  -1.89%  '

$ws.Range("D38").Value = '0.01973'
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("D39").Value = '0.05302'
$ws.Range("E39").Value = '  -1.00%  '

$ws.Range("D40").Value = '2.876'
$ws.Range("E40").Value = '  +1.07%  '

$ws.Range("D41").Value = '0.5173'
$ws.Range("E41").Value = '  -0.52%  '

$ws.Range("D42").Value = '6.983'
$ws.Range("E42").Value = '  +1.83%  '

$ws.Range("D43").Value = '0.1675'
$ws.Range("E43").Value = '  -1.25%  '

$ws.Range("D44").Value = '8.699'
$ws.Range("E44").Value = '  -0.15%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '110.09'
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.76'
$ws.Range("E46").Value = '  +0.86%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.4746'
$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.704'
$ws.Range("E48").Value = '  -1.45%  '

$ws.Range("D49").Value = '1.019'
$ws.Range("E49").Value = '  -1.76%  '

$ws.Range("D50").Value = '0.06473'
$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("D51").Value = '1.849'
$ws.Range("E51").Value = '  -0.40%  '

# Restore the default style on column D so only the displayed value changed
# (the text number-format override is removed again).
$dRange.Style = "Normal"
